# Update column F (dSF) values per the repull/push of data + mean calc
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -4
$ws.Range("F5").Value = 0
$ws.Range("F7").Value = -3
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = -3
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = -3
$ws.Range("F12").Value = -4
$ws.Range("F13").Value = 5
$ws.Range("F14").Value = -2
$ws.Range("F16").Value = 4
$ws.Range("F17").Value = 2
$ws.Range("F19").Value = 5
$ws.Range("F20").Value = -3
$ws.Range("F21").Value = -1
$ws.Range("F22").Value = -1
$ws.Range("F23").Value = -1
$ws.Range("F24").Value = 13
$ws.Range("F25").Value = -6
$ws.Range("F26").Value = -8
$ws.Range("F30").Value = -1
$ws.Range("F34").Value = 0
$ws.Range("F35").Value = -1
$ws.Range("F36").Value = -1
$ws.Range("F37").Value = 1
